$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update query timestamps on the "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:35:43.805831"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:43.805839"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:43.805842"

# --- Add the new "metadata" sheet right after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Match the page margins used elsewhere in the workbook (inches -> points)
$metaSheet.PageSetup.LeftMargin = 54
$metaSheet.PageSetup.RightMargin = 54
$metaSheet.PageSetup.TopMargin = 72
$metaSheet.PageSetup.BottomMargin = 72
$metaSheet.PageSetup.HeaderMargin = 36
$metaSheet.PageSetup.FooterMargin = 36

# Header row (bold, bordered, centered - mirrors the "data" sheet header style)
$header = $metaSheet.Range("B1:G1")
$header.Value = "x"
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# Data row
$a2 = $metaSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

$metaSheet.Range("B2").Value = "Sick sinus syndrome"
$metaSheet.Range("C2").Value = 175

# data_version must stay textual "1.0" (not numeric 1)
$d2 = $metaSheet.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "1.0"
$d2.NumberFormat = "General"

$metaSheet.Range("E2").Value = "2021-08-01T06:00:24.312906Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:43.802108"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/175/?format=json"

Write-Output "metadata sheet added; sheets=$($wb.Worksheets.Count)"
